# Biofuel calculation for Transport
#
# The "Drop-in Diesel" fuel row gets renamed/re-pointed to become the
# "Drop-In Diesel" entry that shares Diesel's color/symbol, and a brand new
# "Drop-In Jet" row is appended that shares Jet Fuel's color/symbol.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 ("Drop-in Diesel") becomes "Drop-In Diesel" and now takes on the
# same color (#58595b) and symbol (diamond) as the "Diesel" row (row 4).
$ws.Range("A13").Value = "Drop-In Diesel"
$ws.Range("B13").Value = "#58595b"
$ws.Range("C13").Value = "diamond"

# New row 31: "Drop-In Jet", sharing "Jet Fuel"'s (row 20) color (#f57e20)
# and symbol (circle).
$ws.Range("A31").Value = "Drop-In Jet"
$ws.Range("B31").Value = "#f57e20"
$ws.Range("C31").Value = "circle"

# Leave the selection on the newly added row, matching the edited file.
$ws.Range("B31:C31").Select()
